$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the old row 6 ("Confirm Order"), pushing it (and the
# trailing blank row) down by one.
$ws.Rows("6:6").Insert()

# Copy the formatting of row 5 into the freshly-inserted row 6 so the new
# row's cell styles line up with the rest of the table.
$ws.Range("A5:D5").Copy($ws.Range("A6:D6"))

# Row 6 content: new "Fiyatin hesaplanmasi" test case.
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Fiyatın hesaplanması"
$ws.Range("C6").Value = "Fiyat ürünü yapmak için gereken  süt, kahve ve su oranlarının, veritabanında belirtilen unit price ile çarpılması durumunda elde edilmesi gerekir. Örneğin: bir kahve yapmak için gerekli olan süt miktarı 2.5 unit, kahve miktarı 1.5 unit, su miktarı 1 unit'tir. Veritabanında ise belirtilen unit price'lar süt için 2 olsun, kahve için 3, su miktarı için 1.5. Bizim elde edeceğimiz fiyat 11 dir. Eğer ekstradan süt eklenirse elde edeceğimiz fiyat 13, kahve eklenirse 14, su eklenirse 12.5 olmalıdır."
$ws.Rows("6:6").RowHeight = 59.4

# The old row 6 ("Confirm Order") now lives at row 7, with its original
# content intact; only its "Test No" cell becomes the literal text "s".
$ws.Range("A7").Value = "s"

$ws.Range("A7").Select()
